# Update the "time_taken" column (F) in the "data" sheet with refreshed
# timestamps, and add a new "metadata" sheet (placed after "data") that
# records details about the PanelApp query that produced this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------
# 1. Refresh the time_taken values on rows 2-43 of the "data" sheet.
# ---------------------------------------------------------------------
$timeValues = @(
  "2021-10-05 14:21:39.000148",
  "2021-10-05 14:21:39.000158",
  "2021-10-05 14:21:39.000162",
  "2021-10-05 14:21:39.000164",
  "2021-10-05 14:21:39.000168",
  "2021-10-05 14:21:39.000170",
  "2021-10-05 14:21:39.000173",
  "2021-10-05 14:21:39.000176",
  "2021-10-05 14:21:39.000179",
  "2021-10-05 14:21:39.000182",
  "2021-10-05 14:21:39.000184",
  "2021-10-05 14:21:39.000187",
  "2021-10-05 14:21:39.000190",
  "2021-10-05 14:21:39.000193",
  "2021-10-05 14:21:39.000196",
  "2021-10-05 14:21:39.000198",
  "2021-10-05 14:21:39.000201",
  "2021-10-05 14:21:39.000204",
  "2021-10-05 14:21:39.000207",
  "2021-10-05 14:21:39.000210",
  "2021-10-05 14:21:39.000213",
  "2021-10-05 14:21:39.000216",
  "2021-10-05 14:21:39.000219",
  "2021-10-05 14:21:39.000221",
  "2021-10-05 14:21:39.000224",
  "2021-10-05 14:21:39.000228",
  "2021-10-05 14:21:39.000230",
  "2021-10-05 14:21:39.000233",
  "2021-10-05 14:21:39.000236",
  "2021-10-05 14:21:39.000239",
  "2021-10-05 14:21:39.000241",
  "2021-10-05 14:21:39.000244",
  "2021-10-05 14:21:39.000248",
  "2021-10-05 14:21:39.000251",
  "2021-10-05 14:21:39.000254",
  "2021-10-05 14:21:39.000256",
  "2021-10-05 14:21:39.000259",
  "2021-10-05 14:21:39.000262",
  "2021-10-05 14:21:39.000265",
  "2021-10-05 14:21:39.000268",
  "2021-10-05 14:21:39.000271",
  "2021-10-05 14:21:39.000274"
)

for ($i = 0; $i -lt $timeValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeValues[$i]
}

# ---------------------------------------------------------------------
# 2. Add the new "metadata" worksheet right after "data".
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"

# Copy the header cell formatting (bold font + border + centered
# alignment, style index 1) from the "data" sheet's header row so the
# new header row looks the same.
$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Copy the formatting of the index column (A2, style index 1) from the
# "data" sheet for the new sheet's index cell.
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Multi-organ autoimmune diabetes"
$meta.Range("C2").Value = 87

# data_version ("1.8") must be stored as text, not as a number.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.8"

$meta.Range("E2").Value = "2020-01-21T16:56:08.048654Z"
$meta.Range("F2").Value = "2021-10-05 14:21:38.996908"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/87/?format=json"

# Make sure "data" stays the active sheet, matching the original workbook.
$ws.Activate()
